# Fin de journée 09/05/2019
# Update "Planning effectif" sheet hours and refresh both sheets' selections.

$wb = $excel.ActiveWorkbook

$wsPrev = $wb.Worksheets.Item("Planning prévisionnel")
$wsEff  = $wb.Worksheets.Item("Planning effectif")

# --- Planning effectif: fill in worked-hours durations (column E) ---
$wsEff.Range("E12").Value = 0.041666666666666664
$wsEff.Range("E13").Value = 0.041666666666666664
$wsEff.Range("E14").Value = 0.041666666666666664
$wsEff.Range("E15").Value = 0.041666666666666664
$wsEff.Range("E17").Value = 0.03125
$wsEff.Range("E18").Value = 0.020833333333333332
$wsEff.Range("E19").Value = 0.041666666666666664
$wsEff.Range("E53").Value = 0.041666666666666664
$wsEff.Range("E56").Value = 0.03125

$excel.Calculate()

# --- Restore selections as saved by each sheet ---
[void]$wsPrev.Range("A20").Select()
[void]$wsEff.Range("F34").Select()
